$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "A slide" -> split "A " into "A" + " " ---
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Characters(1, 1).Text = "A"

# --- TextBox 3: "Followed by a picture" -> split each word/space run ---
$capShape = $s.Shapes.Item(4)
$capRange = $capShape.TextFrame.TextRange
$capRange.Characters(1, 8).Text = "Followed"
$capRange.Characters(10, 2).Text = "by"
$capRange.Characters(13, 1).Text = "a"
